# aktualizace oponentury 3. iterace SIP
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole "BDM" block: row 22 (A=BDM, B=BDM, D="chybi zapocet")
# plus the blank spacer row right below it (old row 23). Deleting both
# rows shifts every following row up by two, matching the target layout
# (old row 24 "Use Case" -> new row 22, old row 30 -> new row 28, and
# the sheet dimension shrinks from J30 to J28).
$ws.Range("A22:A23").EntireRow.Delete()

$ws.Range("A6").Select()
